$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unprotect so the cells can be edited.
$ws.Unprotect()

# Update the confidential disclaimer text (date 2021-03-18 -> 2021-03-19)
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# Keep the original row height (Excel may auto-fit the row after the text edit).
$ws.Rows.Item(7).RowHeight = 15

# Update the numeric values in columns D and E for rows 2-4
$ws.Range("D2").Value = 0.8425326828693127
$ws.Range("E2").Value = 0.002755182366833031

$ws.Range("D3").Value = 0.1574673171306873
$ws.Range("E3").Value = 0.0123549232497191

$ws.Range("E4").Value = 0.004266827808810714

# Restore sheet protection to match the original protected state.
$ws.Protect($null, $true, $true, $true)
